# Locate the "Density analysis" FIELD/YEAR table (4 columns: FIELD/YEAR,
# COUNT, AREA (HA), DENSITY). It is the 16th table in the document.
$d = $word.ActiveDocument
$t = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Columns.Count -eq 4 -and $candidate.Cell(1, 1).Range.Text.StartsWith("FIELD/YEAR") -and $candidate.Cell(1, 2).Range.Text.StartsWith("COUNT")) {
        $t = $candidate
    }
}

# Center every paragraph in every existing cell (header + all data rows).
# Using Paragraphs.Item(1).Alignment (rather than Range.ParagraphFormat.Alignment)
# preserves the existing "Compact" paragraph style on each cell.
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Paragraphs.Item(1).Alignment = 1
    }
}

# Append the new "OKSIR/all years" summary row.
$t.Rows.Add() | Out-Null
$newRowIdx = $t.Rows.Count

$t.Cell($newRowIdx, 1).Range.Text = "OKSIR/all years"
$t.Cell($newRowIdx, 2).Range.Text = "123 - 408"
$t.Cell($newRowIdx, 3).Range.Text = "1383 - 2679"
$t.Cell($newRowIdx, 4).Range.Text = "0.082 - 0.152"

for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $t.Cell($newRowIdx, $c).Range.Paragraphs.Item(1).Alignment = 1
}
